$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# The report repeats the same 2-row block (a data row followed by a
# "Daily Total" row) five times - once per week. Originally only the
# first week's block (rows 5-6) carried the thin-border formatting;
# the other four blocks' value cells had no explicit style. Bring
# every block's formatting into line with week 1's, then add a
# 2-decimal number format across all of the value cells.
# -----------------------------------------------------------------
$dataRows  = @(5, 11, 17, 22, 28)
$totalRows = @(6, 12, 18, 23, 29)

# 1) Task-name cell (column A) of each data row -> thin border only
#    (same look as A5). Donor A5 already has that exact style.
foreach ($r in $dataRows) {
    $ws.Range("A5").Copy() | Out-Null
    $ws.Range("A${r}").PasteSpecial(-4122) | Out-Null
}

# 2) Weekday value cells (columns B:H) of each data row -> thin
#    border (same donor, A5).
foreach ($r in $dataRows) {
    $ws.Range("A5").Copy() | Out-Null
    $ws.Range("B${r}:H${r}").PasteSpecial(-4122) | Out-Null
}

# 3) Weekly-total cell (column I) of each data row -> Courier font +
#    thin border (same look as I5).
foreach ($r in $dataRows) {
    $ws.Range("I5").Copy() | Out-Null
    $ws.Range("I${r}").PasteSpecial(-4122) | Out-Null
}

# 4) "Daily Total" label cell (column A) of each totals row -> bold +
#    thin border (same look as A6).
foreach ($r in $totalRows) {
    $ws.Range("A6").Copy() | Out-Null
    $ws.Range("A${r}").PasteSpecial(-4122) | Out-Null
}

# 5) Totals-row value cells (columns B:I) -> thin border (same donor,
#    A5).
foreach ($r in $totalRows) {
    $ws.Range("A5").Copy() | Out-Null
    $ws.Range("B${r}:I${r}").PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# -----------------------------------------------------------------
# 6) Fix the decimal export: every weekday/weekly-total value cell
#    (data row + totals row, columns B:I) gets a 2-decimal number
#    format.
# -----------------------------------------------------------------
foreach ($r in $dataRows) {
    $ws.Range("B${r}:I${r}").NumberFormat = "0.00"
}
foreach ($r in $totalRows) {
    $ws.Range("B${r}:I${r}").NumberFormat = "0.00"
}

# -----------------------------------------------------------------
# 7) Move the active selection from the title row to B5.
# -----------------------------------------------------------------
$ws.Range("B5").Select()
